$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
Write-Host $ws.Name
